$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
# A "Code" column is inserted after "Name"; "Content" and "Objective" columns
# are dropped; everything else shifts one slot to the left to close the gap.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Code"
$ws.Range("C1").Value = "Course"
$ws.Range("D1").Value = "Teacher"
$ws.Range("E1").Value = "Schedule"
$ws.Range("F1").Value = "Room"
$ws.Range("G1").Value = "AcademicYear"
$ws.Range("H1").Value = "Semester"

# --- Sample data row (row 2) ----------------------------------------------
$ws.Range("A2").Value = "I-ABSCS"
$ws.Range("B2").Value = "I100"
$ws.Range("C2").Value = "SDF"
$ws.Range("D2").Value = "Teacher"
$ws.Range("E2").Value = "MWF"

# "101" looks numeric, so force it to stay text (matching the original
# shared-string "Room" value) before resetting the style back to the
# sheet's default so no stray style index lingers on the cell.
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "101"
$ws.Range("F2").Style = "Normal"

$ws.Range("G2").Value = "2024-2025"
$ws.Range("H2").Value = "First"

# --- Drop the now unused 9th column (old "Semester" duplicate slot) -------
$ws.Columns.Item(9).Delete()
